# Insert a new weekly record at row 201 for
# "Hortaliza, Terminal La Palmera de La Serena - Zanahoria", shifting the
# existing rows 201-250 down to 202-251 (dimension grows from A1:R250 to
# A1:R251).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 201, pushing rows 201..250
# down to 202..251 (carrying their original formatting/values with them).
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new record.
$ws.Cells.Item(201, 1).Value  = 8
$ws.Cells.Item(201, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(201, 3).Value  = "Coquimbo"
$ws.Cells.Item(201, 4).Value  = 44551
$ws.Cells.Item(201, 5).Value  = 4
$ws.Cells.Item(201, 6).Value  = 100114013
$ws.Cells.Item(201, 7).Value  = "Zanahoria"
$ws.Cells.Item(201, 8).Value  = "Sin especificar"
$ws.Cells.Item(201, 9).Value  = "Primera"
$ws.Cells.Item(201, 10).Value = 540
$ws.Cells.Item(201, 11).Value = 6000
$ws.Cells.Item(201, 12).Value = 6500
$ws.Cells.Item(201, 13).Value = 6250
$ws.Cells.Item(201, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(201, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(201, 16).Value = 312
$ws.Cells.Item(201, 17).Value = 20
$ws.Cells.Item(201, 18).Value = "Hortaliza"
